$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numeric-looking plate/ticket ids ("1345", "1111") that
# must stay text (matches the inlineStr string cells in the source data,
# same as row 1's "123"). Force text format first so Excel doesn't
# auto-coerce the values to numbers.
$ws.Range("A2:A3").NumberFormat = "@"

$ws.Range("A2").Value = "1345"
$ws.Range("B2").Value = "car"
$ws.Range("C2").Value = "blue"
$ws.Range("D2").Value = "top"
$ws.Range("E2").Value = "aditya"
$ws.Range("F2").Value = 1768836913.170677
$ws.Range("G2").Value = "1345-2566"

$ws.Range("A3").Value = "1111"
$ws.Range("B3").Value = "car"
$ws.Range("C3").Value = "blue"
$ws.Range("D3").Value = "top"
$ws.Range("E3").Value = "aditya"
$ws.Range("F3").Value = 1768837051.922675
$ws.Range("G3").Value = "1111-1933"
